$wb = $excel.ActiveWorkbook

# zh-cn sheet: update "Correspond Handoff Datetime" (E4) and
# "Correspond Handback DateTime" (H4) for the ffb65260-... row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-12 04:35:56"
$wsZh.Range("H4").Value = "2016-03-12 04:36:14"

# de-de sheet: same two columns for the ffb65260-... row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-12 04:35:59"
$wsDe.Range("H4").Value = "2016-03-12 04:36:19"
